$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: tiny float precision change in single column
$ws.Range("B2").Value = 0.828251637998927

# Row 3 - RandomForestRegressor: updated metrics
$ws.Range("B3").Value = 0.9990288083800415
$ws.Range("C3").Value = 0.998972286036827
$ws.Range("D3").Value = 0.9783891396668984

# Row 4 - model renamed from GradientBoostingRegressor to DecisionTreeRegressor, metrics updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9984841755306016
$ws.Range("C4").Value = 0.9982344228564468
$ws.Range("D4").Value = 0.9914785285293156

# Row 5 - model renamed from AdaBoostRegressor to MLPRegressor, metrics updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9996088301902594
$ws.Range("C5").Value = 0.9995230757379915
$ws.Range("D5").Value = 0.9992759479253376
